# The diff moves the "_GoBack" bookmark (id 0) from the paragraph that
# precedes the last inline picture to the very first paragraph of the
# document (the "Exporting Blender to Three.js" Heading1 paragraph),
# landing right after <w:pPr> and before the first run.

$d = $word.ActiveDocument

# 1) Remove the bookmark from its current location.
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# 2) Re-create it, collapsed, at the very start of the document.
#
# Inserting a zero-length bookmark exactly at document position 0 tends
# to straddle the following paragraph boundary, so we insert a throw-away
# character at position 0, anchor the new bookmark right after it (still
# effectively "at the start" once the character is removed), and then
# delete that helper character again. The bookmark stays collapsed at
# the very beginning of the paragraph.
$d.Range(0, 0).InsertBefore("Z")
$startRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $startRange)
$d.Range(0, 1).Delete()
